# Adding effort on Business analysis for school project
#
# Adds a new effort-log row (row 4) to the "School" sheet:
#   Date: 24-Oct-2015, Name: Sudarshan Acharya, Category: Business Analysis,
#   Man Hours: 2
# Selecting/activating the School sheet makes it the workbook's active tab
# (moving the "active" state off the "Team Member" sheet), and the
# downstream Summery/Cover formulas recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("School")

# Bring School to the front (mirrors the tabSelected/activeTab move in the diff)
$ws.Activate()

$ws.Range("A4").Value = 42301
$ws.Range("B4").Value = "Sudarshan Acharya"
$ws.Range("C4").Value = "Business Analysis"
$ws.Range("D4").Value = 2

# Leave the selection on I4, matching the committed selection state
$ws.Range("I4").Select()
